# This script reproduces the OOXML diff for unitTest_base_macro3.xlsx:
#   - a new "text" command category is introduced on the hidden "#system" sheet.
#     It is inserted alphabetically into the "target" list (column A) and gets
#     its own data column (previously unused column Y), which pushes the
#     existing "web".."xml" categories (and their command lists) one column to
#     the right (Y->Z, Z->AA, AA->AB, AB->AC, AC->AD, AD->AE).
#   - a new "base" command, outputToCloud(resource), is inserted alphabetically
#     into the "base" command list (column E), between "macro(...)" and
#     "prependText(...)".
#   - all named ranges that describe these list boundaries are updated to match
#     the new layout, and a brand new named range "text" is added.
#
# Because named ranges on this sheet are plain cell references (not Excel
# Tables), inserting real rows/columns does not keep them in sync
# automatically in this environment, so the target layout is produced by
# writing the final cell values directly and then re-pointing the named
# ranges to their new extents.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

$ws.Cells.Item(25, "A").Value = "text"
$ws.Cells.Item(26, "A").Value = "web"
$ws.Cells.Item(27, "A").Value = "webalert"
$ws.Cells.Item(28, "A").Value = "webcookie"
$ws.Cells.Item(29, "A").Value = "ws"
$ws.Cells.Item(30, "A").Value = "ws.async"
$ws.Cells.Item(31, "A").Value = "xml"
$ws.Cells.Item(22, "E").Value = "outputToCloud(resource)"
$ws.Cells.Item(23, "E").Value = "prependText(var,prependWith)"
$ws.Cells.Item(24, "E").Value = "repeatUntil(steps,maxWaitMs)"
$ws.Cells.Item(25, "E").Value = "save(var,value)"
$ws.Cells.Item(26, "E").Value = "saveCount(text,regex,saveVar)"
$ws.Cells.Item(27, "E").Value = "saveMatches(text,regex,saveVar)"
$ws.Cells.Item(28, "E").Value = "saveReplace(text,regex,replace,saveVar)"
$ws.Cells.Item(29, "E").Value = "saveVariablesByPrefix(var,prefix)"
$ws.Cells.Item(30, "E").Value = "saveVariablesByRegex(var,regex)"
$ws.Cells.Item(31, "E").Value = "section(steps)"
$ws.Cells.Item(32, "E").Value = "split(text,delim,saveVar)"
$ws.Cells.Item(33, "E").Value = "startRecording()"
$ws.Cells.Item(34, "E").Value = "stopRecording()"
$ws.Cells.Item(35, "E").Value = "substringAfter(text,delim,saveVar)"
$ws.Cells.Item(36, "E").Value = "substringBefore(text,delim,saveVar)"
$ws.Cells.Item(37, "E").Value = "substringBetween(text,start,end,saveVar)"
$ws.Cells.Item(38, "E").Value = "verbose(text)"
$ws.Cells.Item(39, "E").Value = "waitFor(waitMs)"
$ws.Cells.Item(1, "Y").Value = "text"
$ws.Cells.Item(2, "Y").Value = "spellCheck(var,profile,text)"
$ws.Cells.Item(3, "Y").ClearContents()
$ws.Cells.Item(4, "Y").ClearContents()
$ws.Cells.Item(5, "Y").ClearContents()
$ws.Cells.Item(6, "Y").ClearContents()
$ws.Cells.Item(7, "Y").ClearContents()
$ws.Cells.Item(8, "Y").ClearContents()
$ws.Cells.Item(9, "Y").ClearContents()
$ws.Cells.Item(10, "Y").ClearContents()
$ws.Cells.Item(11, "Y").ClearContents()
$ws.Cells.Item(12, "Y").ClearContents()
$ws.Cells.Item(13, "Y").ClearContents()
$ws.Cells.Item(14, "Y").ClearContents()
$ws.Cells.Item(15, "Y").ClearContents()
$ws.Cells.Item(16, "Y").ClearContents()
$ws.Cells.Item(17, "Y").ClearContents()
$ws.Cells.Item(18, "Y").ClearContents()
$ws.Cells.Item(19, "Y").ClearContents()
$ws.Cells.Item(20, "Y").ClearContents()
$ws.Cells.Item(21, "Y").ClearContents()
$ws.Cells.Item(22, "Y").ClearContents()
$ws.Cells.Item(23, "Y").ClearContents()
$ws.Cells.Item(24, "Y").ClearContents()
$ws.Cells.Item(25, "Y").ClearContents()
$ws.Cells.Item(26, "Y").ClearContents()
$ws.Cells.Item(27, "Y").ClearContents()
$ws.Cells.Item(28, "Y").ClearContents()
$ws.Cells.Item(29, "Y").ClearContents()
$ws.Cells.Item(30, "Y").ClearContents()
$ws.Cells.Item(31, "Y").ClearContents()
$ws.Cells.Item(32, "Y").ClearContents()
$ws.Cells.Item(33, "Y").ClearContents()
$ws.Cells.Item(34, "Y").ClearContents()
$ws.Cells.Item(35, "Y").ClearContents()
$ws.Cells.Item(36, "Y").ClearContents()
$ws.Cells.Item(37, "Y").ClearContents()
$ws.Cells.Item(38, "Y").ClearContents()
$ws.Cells.Item(39, "Y").ClearContents()
$ws.Cells.Item(40, "Y").ClearContents()
$ws.Cells.Item(41, "Y").ClearContents()
$ws.Cells.Item(42, "Y").ClearContents()
$ws.Cells.Item(43, "Y").ClearContents()
$ws.Cells.Item(44, "Y").ClearContents()
$ws.Cells.Item(45, "Y").ClearContents()
$ws.Cells.Item(46, "Y").ClearContents()
$ws.Cells.Item(47, "Y").ClearContents()
$ws.Cells.Item(48, "Y").ClearContents()
$ws.Cells.Item(49, "Y").ClearContents()
$ws.Cells.Item(50, "Y").ClearContents()
$ws.Cells.Item(51, "Y").ClearContents()
$ws.Cells.Item(52, "Y").ClearContents()
$ws.Cells.Item(53, "Y").ClearContents()
$ws.Cells.Item(54, "Y").ClearContents()
$ws.Cells.Item(55, "Y").ClearContents()
$ws.Cells.Item(56, "Y").ClearContents()
$ws.Cells.Item(57, "Y").ClearContents()
$ws.Cells.Item(58, "Y").ClearContents()
$ws.Cells.Item(59, "Y").ClearContents()
$ws.Cells.Item(60, "Y").ClearContents()
$ws.Cells.Item(61, "Y").ClearContents()
$ws.Cells.Item(62, "Y").ClearContents()
$ws.Cells.Item(63, "Y").ClearContents()
$ws.Cells.Item(64, "Y").ClearContents()
$ws.Cells.Item(65, "Y").ClearContents()
$ws.Cells.Item(66, "Y").ClearContents()
$ws.Cells.Item(67, "Y").ClearContents()
$ws.Cells.Item(68, "Y").ClearContents()
$ws.Cells.Item(69, "Y").ClearContents()
$ws.Cells.Item(70, "Y").ClearContents()
$ws.Cells.Item(71, "Y").ClearContents()
$ws.Cells.Item(72, "Y").ClearContents()
$ws.Cells.Item(73, "Y").ClearContents()
$ws.Cells.Item(74, "Y").ClearContents()
$ws.Cells.Item(75, "Y").ClearContents()
$ws.Cells.Item(76, "Y").ClearContents()
$ws.Cells.Item(77, "Y").ClearContents()
$ws.Cells.Item(78, "Y").ClearContents()
$ws.Cells.Item(79, "Y").ClearContents()
$ws.Cells.Item(80, "Y").ClearContents()
$ws.Cells.Item(81, "Y").ClearContents()
$ws.Cells.Item(82, "Y").ClearContents()
$ws.Cells.Item(83, "Y").ClearContents()
$ws.Cells.Item(84, "Y").ClearContents()
$ws.Cells.Item(85, "Y").ClearContents()
$ws.Cells.Item(86, "Y").ClearContents()
$ws.Cells.Item(87, "Y").ClearContents()
$ws.Cells.Item(88, "Y").ClearContents()
$ws.Cells.Item(89, "Y").ClearContents()
$ws.Cells.Item(90, "Y").ClearContents()
$ws.Cells.Item(91, "Y").ClearContents()
$ws.Cells.Item(92, "Y").ClearContents()
$ws.Cells.Item(93, "Y").ClearContents()
$ws.Cells.Item(94, "Y").ClearContents()
$ws.Cells.Item(95, "Y").ClearContents()
$ws.Cells.Item(96, "Y").ClearContents()
$ws.Cells.Item(97, "Y").ClearContents()
$ws.Cells.Item(98, "Y").ClearContents()
$ws.Cells.Item(99, "Y").ClearContents()
$ws.Cells.Item(100, "Y").ClearContents()
$ws.Cells.Item(101, "Y").ClearContents()
$ws.Cells.Item(102, "Y").ClearContents()
$ws.Cells.Item(103, "Y").ClearContents()
$ws.Cells.Item(104, "Y").ClearContents()
$ws.Cells.Item(105, "Y").ClearContents()
$ws.Cells.Item(106, "Y").ClearContents()
$ws.Cells.Item(107, "Y").ClearContents()
$ws.Cells.Item(108, "Y").ClearContents()
$ws.Cells.Item(109, "Y").ClearContents()
$ws.Cells.Item(110, "Y").ClearContents()
$ws.Cells.Item(111, "Y").ClearContents()
$ws.Cells.Item(112, "Y").ClearContents()
$ws.Cells.Item(113, "Y").ClearContents()
$ws.Cells.Item(114, "Y").ClearContents()
$ws.Cells.Item(115, "Y").ClearContents()
$ws.Cells.Item(116, "Y").ClearContents()
$ws.Cells.Item(117, "Y").ClearContents()
$ws.Cells.Item(118, "Y").ClearContents()
$ws.Cells.Item(119, "Y").ClearContents()
$ws.Cells.Item(120, "Y").ClearContents()
$ws.Cells.Item(121, "Y").ClearContents()
$ws.Cells.Item(122, "Y").ClearContents()
$ws.Cells.Item(123, "Y").ClearContents()
$ws.Cells.Item(124, "Y").ClearContents()
$ws.Cells.Item(125, "Y").ClearContents()
$ws.Cells.Item(126, "Y").ClearContents()
$ws.Cells.Item(127, "Y").ClearContents()
$ws.Cells.Item(128, "Y").ClearContents()
$ws.Cells.Item(129, "Y").ClearContents()
$ws.Cells.Item(1, "Z").Value = "web"
$ws.Cells.Item(2, "Z").Value = "assertAndClick(locator,label)"
$ws.Cells.Item(3, "Z").Value = "assertAttribute(locator,attrName,value)"
$ws.Cells.Item(4, "Z").Value = "assertAttributeContains(locator,attrName,contains)"
$ws.Cells.Item(5, "Z").Value = "assertAttributeNotContains(locator,attrName,contains)"
$ws.Cells.Item(6, "Z").Value = "assertAttributeNotPresent(locator,attrName)"
$ws.Cells.Item(7, "Z").Value = "assertAttributePresent(locator,attrName)"
$ws.Cells.Item(8, "Z").Value = "assertChecked(locator)"
$ws.Cells.Item(9, "Z").Value = "assertContainCount(locator,text,count)"
$ws.Cells.Item(10, "Z").Value = "assertCssNotPresent(locator,property)"
$ws.Cells.Item(11, "Z").Value = "assertCssPresent(locator,property,value)"
$ws.Cells.Item(12, "Z").Value = "assertElementByAttributes(nameValues)"
$ws.Cells.Item(13, "Z").Value = "assertElementByText(locator,text)"
$ws.Cells.Item(14, "Z").Value = "assertElementCount(locator,count)"
$ws.Cells.Item(15, "Z").Value = "assertElementNotPresent(locator)"
$ws.Cells.Item(16, "Z").Value = "assertElementPresent(locator)"
$ws.Cells.Item(17, "Z").Value = "assertElementsPresent(prefix)"
$ws.Cells.Item(18, "Z").Value = "assertFocus(locator)"
$ws.Cells.Item(19, "Z").Value = "assertFrameCount(count)"
$ws.Cells.Item(20, "Z").Value = "assertFramePresent(frameName)"
$ws.Cells.Item(21, "Z").Value = "assertIECompatMode()"
$ws.Cells.Item(22, "Z").Value = "assertIENativeMode()"
$ws.Cells.Item(23, "Z").Value = "assertLinkByLabel(label)"
$ws.Cells.Item(24, "Z").Value = "assertNotChecked(locator)"
$ws.Cells.Item(25, "Z").Value = "assertNotFocus(locator)"
$ws.Cells.Item(26, "Z").Value = "assertNotText(locator,text)"
$ws.Cells.Item(27, "Z").Value = "assertNotVisible(locator)"
$ws.Cells.Item(28, "Z").Value = "assertOneMatch(locator)"
$ws.Cells.Item(29, "Z").Value = "assertScrollbarHNotPresent(locator)"
$ws.Cells.Item(30, "Z").Value = "assertScrollbarHPresent(locator)"
$ws.Cells.Item(31, "Z").Value = "assertScrollbarVNotPresent(locator)"
$ws.Cells.Item(32, "Z").Value = "assertScrollbarVPresent(locator)"
$ws.Cells.Item(33, "Z").Value = "assertTable(locator,row,column,text)"
$ws.Cells.Item(34, "Z").Value = "assertText(locator,text)"
$ws.Cells.Item(35, "Z").Value = "assertTextContains(locator,text)"
$ws.Cells.Item(36, "Z").Value = "assertTextCount(locator,text,count)"
$ws.Cells.Item(37, "Z").Value = "assertTextList(locator,list,ignoreOrder)"
$ws.Cells.Item(38, "Z").Value = "assertTextMatches(text,minMatch,scrollTo)"
$ws.Cells.Item(39, "Z").Value = "assertTextNotContains(locator,text)"
$ws.Cells.Item(40, "Z").Value = "assertTextNotPresent(text)"
$ws.Cells.Item(41, "Z").Value = "assertTextOrder(locator,descending)"
$ws.Cells.Item(42, "Z").Value = "assertTextPresent(text)"
$ws.Cells.Item(43, "Z").Value = "assertTitle(text)"
$ws.Cells.Item(44, "Z").Value = "assertValue(locator,value)"
$ws.Cells.Item(45, "Z").Value = "assertValueOrder(locator,descending)"
$ws.Cells.Item(46, "Z").Value = "assertVisible(locator)"
$ws.Cells.Item(47, "Z").Value = "checkAll(locator)"
$ws.Cells.Item(48, "Z").Value = "clearLocalStorage()"
$ws.Cells.Item(49, "Z").Value = "click(locator)"
$ws.Cells.Item(50, "Z").Value = "clickAll(locator)"
$ws.Cells.Item(51, "Z").Value = "clickAndWait(locator,waitMs)"
$ws.Cells.Item(52, "Z").Value = "clickByLabel(label)"
$ws.Cells.Item(53, "Z").Value = "clickByLabelAndWait(label,waitMs)"
$ws.Cells.Item(54, "Z").Value = "clickOffset(locator,x,y)"
$ws.Cells.Item(55, "Z").Value = "clickWithKeys(locator,keys)"
$ws.Cells.Item(56, "Z").Value = "close()"
$ws.Cells.Item(57, "Z").Value = "closeAll()"
$ws.Cells.Item(58, "Z").Value = "deselect(locator,text)"
$ws.Cells.Item(59, "Z").Value = "deselectMulti(locator,array)"
$ws.Cells.Item(60, "Z").Value = "dismissInvalidCert()"
$ws.Cells.Item(61, "Z").Value = "dismissInvalidCertPopup()"
$ws.Cells.Item(62, "Z").Value = "doubleClick(locator)"
$ws.Cells.Item(63, "Z").Value = "doubleClickAndWait(locator,waitMs)"
$ws.Cells.Item(64, "Z").Value = "doubleClickByLabel(label)"
$ws.Cells.Item(65, "Z").Value = "doubleClickByLabelAndWait(label,waitMs)"
$ws.Cells.Item(66, "Z").Value = "dragAndDrop(fromLocator,toLocator)"
$ws.Cells.Item(67, "Z").Value = "dragTo(fromLocator,xOffset,yOffset)"
$ws.Cells.Item(68, "Z").Value = "editLocalStorage(key,value)"
$ws.Cells.Item(69, "Z").Value = "executeScript(var,script)"
$ws.Cells.Item(70, "Z").Value = "focus(locator)"
$ws.Cells.Item(71, "Z").Value = "goBack()"
$ws.Cells.Item(72, "Z").Value = "goBackAndWait()"
$ws.Cells.Item(73, "Z").Value = "maximizeWindow()"
$ws.Cells.Item(74, "Z").Value = "mouseOver(locator)"
$ws.Cells.Item(75, "Z").Value = "open(url)"
$ws.Cells.Item(76, "Z").Value = "openAndWait(url,waitMs)"
$ws.Cells.Item(77, "Z").Value = "openHttpBasic(url,username,password)"
$ws.Cells.Item(78, "Z").Value = "openIgnoreTimeout(url)"
$ws.Cells.Item(79, "Z").Value = "refresh()"
$ws.Cells.Item(80, "Z").Value = "refreshAndWait()"
$ws.Cells.Item(81, "Z").Value = "resizeWindow(width,height)"
$ws.Cells.Item(82, "Z").Value = "rightClick(locator)"
$ws.Cells.Item(83, "Z").Value = "saveAllWindowIds(var)"
$ws.Cells.Item(84, "Z").Value = "saveAllWindowNames(var)"
$ws.Cells.Item(85, "Z").Value = "saveAttribute(var,locator,attrName)"
$ws.Cells.Item(86, "Z").Value = "saveAttributeList(var,locator,attrName)"
$ws.Cells.Item(87, "Z").Value = "saveCount(var,locator)"
$ws.Cells.Item(88, "Z").Value = "saveDivsAsCsv(headers,rows,cells,nextPage,file)"
$ws.Cells.Item(89, "Z").Value = "saveElement(var,locator)"
$ws.Cells.Item(90, "Z").Value = "saveElements(var,locator)"
$ws.Cells.Item(91, "Z").Value = "saveLocalStorage(var,key)"
$ws.Cells.Item(92, "Z").Value = "saveLocation(var)"
$ws.Cells.Item(93, "Z").Value = "savePageAs(var,sessionIdName,url)"
$ws.Cells.Item(94, "Z").Value = "savePageAsFile(sessionIdName,url,file)"
$ws.Cells.Item(95, "Z").Value = "saveTableAsCsv(locator,nextPageLocator,file)"
$ws.Cells.Item(96, "Z").Value = "saveText(var,locator)"
$ws.Cells.Item(97, "Z").Value = "saveTextArray(var,locator)"
$ws.Cells.Item(98, "Z").Value = "saveTextSubstringAfter(var,locator,delim)"
$ws.Cells.Item(99, "Z").Value = "saveTextSubstringBefore(var,locator,delim)"
$ws.Cells.Item(100, "Z").Value = "saveTextSubstringBetween(var,locator,start,end)"
$ws.Cells.Item(101, "Z").Value = "saveValue(var,locator)"
$ws.Cells.Item(102, "Z").Value = "saveValues(var,locator)"
$ws.Cells.Item(103, "Z").Value = "scrollElement(locator,xOffset,yOffset)"
$ws.Cells.Item(104, "Z").Value = "scrollLeft(locator,pixel)"
$ws.Cells.Item(105, "Z").Value = "scrollPage(xOffset,yOffset)"
$ws.Cells.Item(106, "Z").Value = "scrollRight(locator,pixel)"
$ws.Cells.Item(107, "Z").Value = "scrollTo(locator)"
$ws.Cells.Item(108, "Z").Value = "select(locator,text)"
$ws.Cells.Item(109, "Z").Value = "selectFrame(locator)"
$ws.Cells.Item(110, "Z").Value = "selectMulti(locator,array)"
$ws.Cells.Item(111, "Z").Value = "selectMultiOptions(locator)"
$ws.Cells.Item(112, "Z").Value = "selectText(locator)"
$ws.Cells.Item(113, "Z").Value = "selectWindow(winId)"
$ws.Cells.Item(114, "Z").Value = "selectWindowAndWait(winId,waitMs)"
$ws.Cells.Item(115, "Z").Value = "selectWindowByIndex(index)"
$ws.Cells.Item(116, "Z").Value = "selectWindowByIndexAndWait(index,waitMs)"
$ws.Cells.Item(117, "Z").Value = "toggleSelections(locator)"
$ws.Cells.Item(118, "Z").Value = "type(locator,value)"
$ws.Cells.Item(119, "Z").Value = "typeKeys(locator,value)"
$ws.Cells.Item(120, "Z").Value = "uncheckAll(locator)"
$ws.Cells.Item(121, "Z").Value = "unselectAllText()"
$ws.Cells.Item(122, "Z").Value = "upload(fieldLocator,file)"
$ws.Cells.Item(123, "Z").Value = "verifyContainText(locator,text)"
$ws.Cells.Item(124, "Z").Value = "verifyText(locator,text)"
$ws.Cells.Item(125, "Z").Value = "wait(waitMs)"
$ws.Cells.Item(126, "Z").Value = "waitForElementPresent(locator)"
$ws.Cells.Item(127, "Z").Value = "waitForPopUp(winId,waitMs)"
$ws.Cells.Item(128, "Z").Value = "waitForTextPresent(text)"
$ws.Cells.Item(129, "Z").Value = "waitForTitle(text)"
$ws.Cells.Item(1, "AA").Value = "webalert"
$ws.Cells.Item(2, "AA").Value = "accept()"
$ws.Cells.Item(3, "AA").Value = "assertPresent()"
$ws.Cells.Item(4, "AA").Value = "assertText(text,matchBy)"
$ws.Cells.Item(5, "AA").Value = "dismiss()"
$ws.Cells.Item(6, "AA").Value = "replyCancel(text)"
$ws.Cells.Item(7, "AA").Value = "replyOK(text)"
$ws.Cells.Item(8, "AA").Value = "storeText(var)"
$ws.Cells.Item(1, "AB").Value = "webcookie"
$ws.Cells.Item(2, "AB").Value = "assertNotPresent(name)"
$ws.Cells.Item(3, "AB").Value = "assertPresent(name)"
$ws.Cells.Item(4, "AB").Value = "assertValue(name,value)"
$ws.Cells.Item(5, "AB").Value = "delete(name)"
$ws.Cells.Item(6, "AB").Value = "deleteAll()"
$ws.Cells.Item(7, "AB").Value = "save(var,name)"
$ws.Cells.Item(8, "AB").Value = "saveAll(var)"
$ws.Cells.Item(9, "AB").ClearContents()
$ws.Cells.Item(10, "AB").ClearContents()
$ws.Cells.Item(11, "AB").ClearContents()
$ws.Cells.Item(12, "AB").ClearContents()
$ws.Cells.Item(13, "AB").ClearContents()
$ws.Cells.Item(14, "AB").ClearContents()
$ws.Cells.Item(15, "AB").ClearContents()
$ws.Cells.Item(16, "AB").ClearContents()
$ws.Cells.Item(17, "AB").ClearContents()
$ws.Cells.Item(1, "AC").Value = "ws"
$ws.Cells.Item(2, "AC").Value = "assertReturnCode(var,returnCode)"
$ws.Cells.Item(3, "AC").Value = "delete(url,body,var)"
$ws.Cells.Item(4, "AC").Value = "download(url,queryString,saveTo)"
$ws.Cells.Item(5, "AC").Value = "get(url,queryString,var)"
$ws.Cells.Item(6, "AC").Value = "head(url,var)"
$ws.Cells.Item(7, "AC").Value = "header(name,value)"
$ws.Cells.Item(8, "AC").Value = "headerByVar(name,var)"
$ws.Cells.Item(9, "AC").Value = "jwtParse(var,token,key)"
$ws.Cells.Item(10, "AC").Value = "jwtSignHS256(var,payload,key)"
$ws.Cells.Item(11, "AC").Value = "oauth(var,url,auth)"
$ws.Cells.Item(12, "AC").Value = "patch(url,body,var)"
$ws.Cells.Item(13, "AC").Value = "post(url,body,var)"
$ws.Cells.Item(14, "AC").Value = "put(url,body,var)"
$ws.Cells.Item(15, "AC").Value = "saveResponsePayload(var,file,append)"
$ws.Cells.Item(16, "AC").Value = "soap(action,url,payload,var)"
$ws.Cells.Item(17, "AC").Value = "upload(url,body,fileParams,var)"
$ws.Cells.Item(1, "AD").Value = "ws.async"
$ws.Cells.Item(2, "AD").Value = "delete(url,body,output)"
$ws.Cells.Item(3, "AD").Value = "download(url,queryString,saveTo)"
$ws.Cells.Item(4, "AD").Value = "get(url,queryString,output)"
$ws.Cells.Item(5, "AD").Value = "head(url,output)"
$ws.Cells.Item(6, "AD").Value = "patch(url,body,output)"
$ws.Cells.Item(7, "AD").Value = "post(url,body,output)"
$ws.Cells.Item(8, "AD").Value = "put(url,body,output)"
$ws.Cells.Item(9, "AD").ClearContents()
$ws.Cells.Item(10, "AD").ClearContents()
$ws.Cells.Item(11, "AD").ClearContents()
$ws.Cells.Item(12, "AD").ClearContents()
$ws.Cells.Item(13, "AD").ClearContents()
$ws.Cells.Item(14, "AD").ClearContents()
$ws.Cells.Item(15, "AD").ClearContents()
$ws.Cells.Item(16, "AD").ClearContents()
$ws.Cells.Item(17, "AD").ClearContents()
$ws.Cells.Item(18, "AD").ClearContents()
$ws.Cells.Item(19, "AD").ClearContents()
$ws.Cells.Item(20, "AD").ClearContents()
$ws.Cells.Item(21, "AD").ClearContents()
$ws.Cells.Item(22, "AD").ClearContents()
$ws.Cells.Item(23, "AD").ClearContents()
$ws.Cells.Item(24, "AD").ClearContents()
$ws.Cells.Item(25, "AD").ClearContents()
$ws.Cells.Item(26, "AD").ClearContents()
$ws.Cells.Item(27, "AD").ClearContents()
$ws.Cells.Item(1, "AE").Value = "xml"
$ws.Cells.Item(2, "AE").Value = "append(xml,xpath,content,var)"
$ws.Cells.Item(3, "AE").Value = "assertCorrectness(xml,schema)"
$ws.Cells.Item(4, "AE").Value = "assertElementCount(xml,xpath,count)"
$ws.Cells.Item(5, "AE").Value = "assertElementNotPresent(xml,xpath)"
$ws.Cells.Item(6, "AE").Value = "assertElementPresent(xml,xpath)"
$ws.Cells.Item(7, "AE").Value = "assertSoap(wsdl,xml)"
$ws.Cells.Item(8, "AE").Value = "assertSoapFaultCode(expected,xml)"
$ws.Cells.Item(9, "AE").Value = "assertSoapFaultString(expected,xml)"
$ws.Cells.Item(10, "AE").Value = "assertValue(xml,xpath,expected)"
$ws.Cells.Item(11, "AE").Value = "assertValues(xml,xpath,array,exactOrder)"
$ws.Cells.Item(12, "AE").Value = "assertWellformed(xml)"
$ws.Cells.Item(13, "AE").Value = "beautify(xml,var)"
$ws.Cells.Item(14, "AE").Value = "clear(xml,xpath,var)"
$ws.Cells.Item(15, "AE").Value = "delete(xml,xpath,var)"
$ws.Cells.Item(16, "AE").Value = "insertAfter(xml,xpath,content,var)"
$ws.Cells.Item(17, "AE").Value = "insertBefore(xml,xpath,content,var)"
$ws.Cells.Item(18, "AE").Value = "minify(xml,var)"
$ws.Cells.Item(19, "AE").Value = "prepend(xml,xpath,content,var)"
$ws.Cells.Item(20, "AE").Value = "replace(xml,xpath,content,var)"
$ws.Cells.Item(21, "AE").Value = "replaceIn(xml,xpath,content,var)"
$ws.Cells.Item(22, "AE").Value = "storeCount(xml,xpath,var)"
$ws.Cells.Item(23, "AE").Value = "storeSoapFaultCode(var,xml)"
$ws.Cells.Item(24, "AE").Value = "storeSoapFaultDetail(var,xml)"
$ws.Cells.Item(25, "AE").Value = "storeSoapFaultString(var,xml)"
$ws.Cells.Item(26, "AE").Value = "storeValue(xml,xpath,var)"
$ws.Cells.Item(27, "AE").Value = "storeValues(xml,xpath,var)"

# --- Re-point the named ranges (definedNames) to the new layout ---
$wb.Names.Item("base").RefersTo = "='#system'!`$E`$2:`$E`$39"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$31"
$wb.Names.Item("web").RefersTo = "='#system'!`$Z`$2:`$Z`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AC`$2:`$AC`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AD`$2:`$AD`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AE`$2:`$AE`$27"

# --- Add the brand new named range for the "text" category ---
$wb.Names.Add("text", "='#system'!`$Y`$2:`$Y`$2")
